$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The job card had two extra items (rows for items 5 and 6) that are no
# longer needed — remove them, shrinking the used range from A1:K7 to A1:K5.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# The remaining "Total SQ.FT." entries for items 2 and 3 now carry values.
$ws.Range("G3").Value = 7.719
$ws.Range("G4").Value = 11.579
